# Release v0.1.0-beta: Fix validation errors and update canonical URL
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Version: 1.0.0 -> 0.1.0
$ws.Range("B3").Value = "0.1.0"

# Status: active -> draft
$ws.Range("B6").Value = "draft"

# Experimental: (empty) -> "false" (must stay a literal text value, not a
# boolean -- writing it directly would get auto-coerced to a Boolean cell,
# so stage it via an apostrophe-forced text cell and PasteSpecial the value
# back in, which keeps the original cell style intact).
$ws.Range("Z99").Value = "'false"
$ws.Range("Z99").Copy()
$ws.Range("B7").PasteSpecial(-4163)
$ws.Range("Z99").Clear()

# Date: updated timestamp
$ws.Range("B8").Value = "2025-12-26T14:13:58+00:00"

# Description: (empty) -> descriptive text
$ws.Range("B11").Value = "Value set for categorizing nursing problems"
